# Append a new data row (row 89) to each of the 4 worksheets, mirroring
# the structure/format of the existing rows (time, length, ID, actual
# length, checksum + their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rows = @{
    "MID_LFT_#1" = @{
        A = 45875.45994212963
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x20"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 288
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45875.45994212963
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x2C"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 300
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45875.45994212963
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5F"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 95
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45875.45994212963
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 117
        I = 9
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 89

    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $data.A
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
